$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rename the sheets (Sheet1/2/3 -> Assignments/StudentDetails/RegistrationURL)
# ---------------------------------------------------------------------------
$wsAssignments = $wb.Worksheets.Item(1)
$wsStudents    = $wb.Worksheets.Item(2)
$wsRegUrl      = $wb.Worksheets.Item(3)

$wsAssignments.Name = "Assignments"
$wsStudents.Name    = "StudentDetails"
$wsRegUrl.Name      = "RegistrationURL"

# ---------------------------------------------------------------------------
# Populate the "StudentDetails" sheet header row, then the "RegistrationURL"
# header, then the rest of the StudentDetails/RegistrationURL data - this is
# the order the data was actually typed in, and it drives the order new
# strings land in the shared string table.
# ---------------------------------------------------------------------------
$wsStudents.Range("A1").Value = "StudentId"
$wsStudents.Range("B1").Value = "StudentPwd"
$wsStudents.Range("C1").Value = "FirstName"
$wsStudents.Range("D1").Value = "LastName"
$wsStudents.Range("A1:D1").Font.Bold = $true

$wsRegUrl.Range("A1").Value = "registrationURL"
$wsRegUrl.Range("A1").Font.Bold = $true

$wsStudents.Range("B2").Value = "Aa123456"
$wsStudents.Range("D2").Value = "STUDENT 1"

$wsRegUrl.Range("B3").Value = "http://connectqastg.mheducation.com/class/k-ins-chemistry-section1"
$wsRegUrl.Range("B1").Value = "http://connectqastg.mheducation.com/class/-_-_123--fname--section-1-2"

# The StudentId cell is a real mailto: hyperlink (this also creates the
# built-in "Hyperlink" cell style used below).
$wsStudents.Hyperlinks.Add($wsStudents.Range("A2"), "mailto:stgauguststudent1@gmail.com", "", "", "stgauguststudent1@gmail.com") | Out-Null

$wsStudents.Range("C2").Value = "STG AUGUST"

# RegistrationURL!B1 only borrows the Hyperlink look (blue/underline) - it is
# not an actual clickable hyperlink.
$wsRegUrl.Range("B1").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column widths (best-fit in the original file) - approximate via ColumnWidth
# ---------------------------------------------------------------------------
$wsStudents.Columns.Item(1).ColumnWidth = 27.736979166666668
$wsStudents.Columns.Item(2).ColumnWidth = 10.877604166666666
$wsStudents.Columns.Item(3).ColumnWidth = 16.451822916666668
$wsStudents.Columns.Item(4).ColumnWidth = 67.16666666666667

$wsRegUrl.Columns.Item(1).ColumnWidth = 13.877604166666666
$wsRegUrl.Columns.Item(2).ColumnWidth = 68.16666666666667

# ---------------------------------------------------------------------------
# Page setup for RegistrationURL (keeps it in sync with the other sheets)
# ---------------------------------------------------------------------------
$wsRegUrl.PageSetup.PaperSize = 9
$wsRegUrl.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selections / active sheet / active cell bookkeeping
# ---------------------------------------------------------------------------
$wsAssignments.Range("B31").Select()
$wsStudents.Range("D20").Select()
$wsRegUrl.Range("B20").Select()

$wsStudents.Activate()
